$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'70.502.61"
$ws.Range("E2").Value = "  +2.02%  "
$ws.Range("D3").Value = "'3.561.18"
$ws.Range("E3").Value = "  +1.25%  "
$ws.Range("D4").Value = "'1.00"
$ws.Range("E4").Value = "  +0.00%  "
$ws.Range("D5").Value = "'613.91"
$ws.Range("E5").Value = "  +6.12%  "
$ws.Range("D6").Value = "'172.52"
$ws.Range("E6").Value = "  +0.63%  "
$ws.Range("D7").Value = "'0.617"
$ws.Range("E7").Value = "  +1.51%  "
$ws.Range("D8").Value = "'3.556.57"
$ws.Range("E8").Value = "  +1.29%  "
$ws.Range("E9").Value = "  -0.03%  "
$ws.Range("E10").Value = "  +4.14%  "
$ws.Range("D11").Value = "'7.23"
$ws.Range("E11").Value = "  +7.36%  "
$ws.Range("E12").Value = "  +0.51%  "
$ws.Range("D13").Value = "'46.57"
$ws.Range("E13").Value = "  -0.81%  "
$ws.Range("D14").Value = "'0.0000277"
$ws.Range("E14").Value = "  +1.31%  "
$ws.Range("D15").Value = "'4.137.98"
$ws.Range("E15").Value = "  +1.41%  "
$ws.Range("D16").Value = "'8.35"
$ws.Range("E16").Value = "  -1.97%  "
$ws.Range("D17").Value = "'615.18"
$ws.Range("E17").Value = "  -0.97%  "
$ws.Range("D18").Value = "'3.562.46"
$ws.Range("E18").Value = "  +1.86%  "
$ws.Range("D19").Value = "'70.629.09"
$ws.Range("E19").Value = "  +2.23%  "
$ws.Range("E20").Value = "  -2.19%  "
$ws.Range("D21").Value = "'17.41"
$ws.Range("E21").Value = "  -0.07%  "
$ws.Range("D22").Value = "'0.880"
$ws.Range("E22").Value = "  -0.42%  "
$ws.Range("E23").Value = "  -15.79%  "
$ws.Range("E24").Value = "  -1.20%  "
$ws.Range("D25").Value = "'96.73"
$ws.Range("E25").Value = "  -0.71%  "
$ws.Range("D26").Value = "'3.83"
$ws.Range("E26").Value = "  +0.85%  "
$ws.Range("E27").Value = "  +0.08%  "
$ws.Range("E28").Value = "  -1.09%  "
$ws.Range("D29").Value = "'33.53"
$ws.Range("E29").Value = "  +2.79%  "
$ws.Range("E30").Value = "  -3.33%  "
$ws.Range("D31").Value = "'8.48"
$ws.Range("E31").Value = "  -0.81%  "
$ws.Range("E32").Value = "  -3.76%  "
$ws.Range("E33").Value = "  -1.36%  "
$ws.Range("D34").Value = "'6.94"
$ws.Range("E34").Value = "  -0.65%  "
$ws.Range("D35").Value = "'572.66"
$ws.Range("E35").Value = "  -9.66%  "
$ws.Range("D36").Value = "'3.63"
$ws.Range("E36").Value = "  +6.33%  "
$ws.Range("E37").Value = "  -1.26%  "
$ws.Range("D38").Value = "'10.81"
$ws.Range("E38").Value = "  +0.51%  "
$ws.Range("D39").Value = "'57.33"
$ws.Range("E39").Value = "  +1.11%  "
$ws.Range("E40").Value = "  +4.90%  "
$ws.Range("E41").Value = "  +0.13%  "
$ws.Range("D43").Value = "'3.382.64"
$ws.Range("E43").Value = "  +0.08%  "
$ws.Range("E44").Value = "  -2.44%  "
$ws.Range("D45").Value = "'32.99"
$ws.Range("E45").Value = "  +0.29%  "
$ws.Range("D46").Value = "'2.95"
$ws.Range("E46").Value = "  +7.20%  "
$ws.Range("D47").Value = "'0.0₃0701"
$ws.Range("E47").Value = "  +1.57%  "
$ws.Range("D48").Value = "'2.61"
$ws.Range("E48").Value = "  +1.84%  "
$ws.Range("E49").Value = "  +0.13%  "
$ws.Range("D50").Value = "'133.71"
$ws.Range("E50").Value = "  +1.13%  "
